# Automatic update of files.
# - Bump the "Förändrad" (Changed) date in column C from 2024-06-06 (45449)
#   to 2024-06-07 (45450) for every data row (rows 2-28).
# - Remove the last data row (row 29, "A 22929-2024") which was the row
#   that previously carried the new "Förändrad" date.
# - Row 28 goes back to the sheet's default row height (its explicit
#   customHeight="1" is cleared) now that it is the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 28 to the new date serial.
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45450
}

# Delete row 29 entirely (shifts nothing up since it's the last row).
$ws.Rows(29).Delete()

# Clear the explicit row height on row 28 so it reverts to the sheet default.
$ws.Rows(28).AutoFit()
